# Applies the two changes captured by the commit:
#   1. Slide 6's table switches from the custom "Table_0" style
#      ({F2ECC8C6-D526-43CA-ABB5-87180C55A4D9}) to the built-in table
#      style {6B575B69-7263-435E-BCC0-9213BF831D1D}.
#   2. The deck's live theme (ppt/theme/theme2.xml, the one actually wired
#      to the slide master) swaps its 12-colour scheme from the "Integral"
#      palette over to the stock "Office Theme" palette (the palette that
#      used to live, unused, in ppt/theme/theme1.xml).

$p = $ppt.ActivePresentation

# --- 1. Retarget the table's style -----------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{6B575B69-7263-435E-BCC0-9213BF831D1D}")
        }
    }
}

# --- 2. Swap the active theme's colour scheme to "Office Theme" ------
# Index order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink. Values are OLE COLORREF (BGR) longs.
$officeThemeColors = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388     # dk2      44546A
    4  = 15132391    # lt2      E7E6E6
    5  = 13998939    # accent1  5B9BD5
    6  = 3243501      # accent2  ED7D31
    7  = 10855845     # accent3  A5A5A5
    8  = 49407         # accent4  FFC000
    9  = 12874308     # accent5  4472C4
    10 = 4697456      # accent6  70AD47
    11 = 12673797     # hlink    0563C1
    12 = 7491477      # folHlink 954F72
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeColors[$i]
}
